$d = $word.ActiveDocument

# Non-breaking space used between "N"/"E" and the coordinate numbers in the
# existing document text.
$nbsp = [char]0x00A0

# ---------------------------------------------------------------------------
# 1. Title paragraph: "höga naturvärden i ..." -> "höga naturvärden och
#    fridlysta arter i ...". Scope the Find to the Title paragraph only,
#    because the same phrase also appears (unchanged) in the next paragraph.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.Find.Execute(
    "höga naturvärden i avverkningsanmälan A 52046-2022 i Södertälje kommun",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "höga naturvärden och fridlysta arter i avverkningsanmälan A 52046-2022 i Södertälje kommun",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Area figure: 8,4 ha -> 14,4 ha
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "och omfattar 8,4 ha.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "och omfattar 14,4 ha.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3. "Nedan beskrivs..." sentence: add "och fridlysta arter"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Nedan beskrivs fynd av naturvårdsarter som gjorts i det avverkningsanmälda området.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Nedan beskrivs fynd av naturvårdsarter och fridlysta arter som gjorts i det avverkningsanmälda området.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Species-found sentence replacement
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "I avverkningsanmälan har följande 1 naturvårdsarter hittats: tjockfotad fingersvamp (S). Arter som är signalarter enligt Skogsstyrelsen har markerats med (S).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I avverkningsanmälan har följande 1 naturvårdsarter hittats: nattviol (§8). För fridlysta arter anges även paragrafen i Artskyddsförordningen som arten är fridlyst enligt.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 5. Map caption coordinates
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Kartans mittpunktskoordinat är N${nbsp}6560698, E${nbsp}647872 i SWEREF 99 TM.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Kartans mittpunktskoordinat är N${nbsp}6560497, E${nbsp}647779 i SWEREF 99 TM.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 6. Insert a new "Fridlysta arter" (Heading1) paragraph right after the
#    figure caption paragraph, and replace the final paragraph (which used to
#    describe "Tjockfotad fingersvamp" with a bold lead-in run) with a single,
#    plain paragraph about the protected species "nattviol".
# ---------------------------------------------------------------------------
$captionPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Figur 1.*SWEREF 99 TM.*") {
        $captionPara = $p
    }
}

$captionPara.Range.InsertParagraphAfter() | Out-Null
$headingPara = $captionPara.Next()
$headingPara.Range.Text = "Fridlysta arter"
$headingPara.Range.ParagraphFormat.Style = "Heading1"

# Insert a fresh (unformatted) paragraph after the new heading, fill it with
# the replacement text, then drop the old "Tjockfotad fingersvamp..."
# paragraph that follows it.
$headingPara.Range.InsertParagraphAfter() | Out-Null
$bodyPara = $headingPara.Next()
$bodyPara.Range.ParagraphFormat.Style = "Normal"
$bodyPara.Range.Text = "Följande fridlysta arter har sina livsmiljöer och växtplatser i den avverkningsanmälda skogen: nattviol (§8)."

$oldSpeciesPara = $bodyPara.Next()
$oldSpeciesPara.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 7. Update the date stamp that appears in the "first page" header.
# ---------------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    foreach ($h in $sec.Headers) {
        if ($h.Exists) {
            $h.Range.Find.Execute(
                "2026-02-13", $true, $false, $false, $false, $false, $true, 1, $false,
                "2026-02-17", 2) | Out-Null
        }
    }
}

